# SSEL-PGCS.docx edit
#  1. Merge the three runs of "Conjunto de componentes con una determinada
#     versión que en forma conjunto permite el funcionamiento de la
#     aplicación." into a single run.
#  2. Drop the stray "_GoBack" bookmark that used to sit next to the
#     h.3dy6vkm bookmark.
#  3. Re-create the "_GoBack" bookmark right after the "Organización de SCM"
#     heading run (this also naturally renumbers the other bookmark ids,
#     same as Word does when a bookmark is removed/added).

$d = $word.ActiveDocument

# --- 1. Merge the "Conjunto de componentes..." runs into one ---------------
$texto = "Conjunto de componentes con una determinada versión que en forma " + `
         "conjunto permite el funcionamiento de la aplicación."
$d.Content.Find.Execute($texto, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $texto, 2) | Out-Null

# --- 2. Remove the old "_GoBack" bookmark -----------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 3. Re-add "_GoBack" right after the body heading "Organización de SCM" -
# The document contains "Organización de SCM" twice: once in the TOC
# (hyperlink) and once as the actual heading further down. We want the
# second (body) occurrence.
$rng = $d.Content
$rng.Find.ClearFormatting()
$hitCount = 0
while ($rng.Find.Execute("Organización de SCM", $true, $false, $false, `
                          $false, $false, $true, 1, $false, "", 0)) {
    $hitCount = $hitCount + 1
    if ($hitCount -eq 2) {
        break
    }
    $rng.Collapse(0)
    $rng.End = $d.Content.End
}

$pt = $rng.End

# Placing a zero-length bookmark exactly on a run/paragraph boundary can
# land in the wrong spot, so insert a throw-away marker character right
# after the heading, wrap the bookmark around that single character, then
# delete the character again - the bookmark collapses back down to a point
# and stays put, right after "Organización de SCM".
$insPoint = $d.Range($pt, $pt)
$insPoint.InsertAfter("X")

$wrap = $d.Range($pt, $pt + 1)
$d.Bookmarks.Add("_GoBack", $wrap)

$marker = $d.Range($pt, $pt + 1)
$marker.Delete()
